# Generate Report for Handback
# Update timestamps (and one priority value) across the three worksheets
# to reflect the new report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-08 23:58:44"
$wsOverview.Range("G3").Value = "2016-11-08 23:58:44"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-11-08 23:58:30"
$wsZhCn.Range("H3").Value = "2016-11-08 23:58:30"
$wsZhCn.Range("K2").Value = "2016-11-08 23:59:22"
$wsZhCn.Range("K3").Value = "2016-11-08 23:59:22"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-11-08 23:58:44"
$wsDeDe.Range("H3").Value = "2016-11-08 23:58:44"
$wsDeDe.Range("K2").Value = "2016-11-08 23:59:40"
$wsDeDe.Range("K3").Value = "2016-11-08 23:59:40"
